$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (row 1): folder / class / regexp / regexp_name ---
$ws.Range("D1").Value2 = "folder"
$ws.Range("E1").Value2 = "class"

# --- Row 48: folder id for existing "Номенклатурная позиция" class row ---
$ws.Range("D48").Value2 = "df0921c1-f46f-e911-8115-817c3f53a992"

# --- Row 47 & 48 share the same class id value ---
$ws.Range("E47").Value2 = "0e1d8277-d859-e911-8115-817c3f53a992"
$ws.Range("E48").Value2 = "0e1d8277-d859-e911-8115-817c3f53a992"

# --- New rows 49 and 50 (additional class/type entries) ---
$ws.Range("A49").Value2 = "c7ec82e0-f360-e911-8115-817c3f53a992"
$ws.Range("B49").Value2 = 2
$ws.Range("C49").Value2 = "Подобъект"

$ws.Range("A50").Value2 = "532d2888-3582-ec11-911c-005056b6948b"
$ws.Range("B50").Value2 = 8
$ws.Range("C50").Value2 = "Номенклатурная позиция"

# --- Regexp used to extract "stamp of doc" / "№ поз. По ГП" for row 50 ---
$ws.Range("F50").Value2 = "^(\S*)\s"

# --- Remaining headers for new columns F/G ---
$ws.Range("F1").Value2 = "regexp"
$ws.Range("G1").Value2 = "regexp_name"

# --- Stamp labels for rows 49/50 ---
$ws.Range("G49").Value2 = "№ поз. По ГП"
$ws.Range("G50").Value2 = "Тип материала"

# --- folder/class ids for new row 50 ---
$ws.Range("D50").Value2 = "52ad3a8a-3382-ec11-911c-005056b6948b"
$ws.Range("E50").Value2 = "708c334d-e78a-ec11-911c-005056b6948b"

# --- folder id for row 47 (keeps the same number format used by A15) ---
$ws.Range("D47").NumberFormat = "0.00E+00"
$ws.Range("D47").Value2 = "7e2318df-0fe9-e911-80cf-9706d383f138"

# --- regexp for row 49, stored as text ---
$ws.Range("F49").NumberFormat = "@"
$ws.Range("F49").Value2 = "^.{0,5}?\.?(\S?\d+(\.\d)?)"

# F50 also uses the text number format
$ws.Range("F50").NumberFormat = "@"

# --- Column widths for the new D:G columns ---
$ws.Columns.Item(4).ColumnWidth = 35.33333333333333
$ws.Columns.Item(5).ColumnWidth = 35.5
$ws.Columns.Item(6).ColumnWidth = 20.333333333333336
$ws.Columns.Item(7).ColumnWidth = 12.0

# --- View state: scrolled down, D46 selected ---
$ws.Range("D46").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 37
$win.ScrollColumn = 1
